$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Shorten supplier names by removing the region suffix.
# "FG EMPREIT. MAO DE OBRA LTDA PU_SUL" -> "FG EMPREIT. MAO DE OBRA LTDA"
# "JAPJ CONSTRUCOES CIVIS LTDA PU_SUDESTE" -> "JAPJ CONSTRUCOES CIVIS"
$ws.Range("A2").Value = "FG EMPREIT. MAO DE OBRA LTDA"
$ws.Range("A3").Value = "JAPJ CONSTRUCOES CIVIS"
$ws.Range("A5").Value = "FG EMPREIT. MAO DE OBRA LTDA"
$ws.Range("A15").Value = "JAPJ CONSTRUCOES CIVIS"

$wb.Save()
